$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above current row 2 (pushes the existing data rows down by one).
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with slugified "machine" codes for each
# column header in row 1. These new codes let two columns be related to each
# other (e.g. "grado-codigo" <-> "grado") to build SKOS hierarchies.
$ws.Range("A2").Value = "grado-codigo"
$ws.Range("B2").Value = "personas"
$ws.Range("C2").Value = "grandes-grupos-codigo"
$ws.Range("D2").Value = "provincia-codigo"
$ws.Range("E2").Value = "comarca-nombre"
$ws.Range("F2").Value = "grado"
$ws.Range("G2").Value = "sexo-codigo"
$ws.Range("H2").Value = "provincia-nombre"
$ws.Range("I2").Value = "sexo"
$ws.Range("J2").Value = "grandes-grupos"
$ws.Range("K2").Value = "comarca-codigo"
